# Fruta / hortaliza, semanal
# Weekly refresh of the "Vega Monumental Concepción - Granada" subset:
# the three data rows (2-4) are rotated by one position -- the newest
# observation (previously in row 4) becomes the first row, and the
# older rows shift down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current values of the columns that actually change
# (D, M, N, O, P, R, S) for rows 2-4 before overwriting anything.
$cols = @("D", "M", "N", "O", "P", "R", "S")

$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range($col + "2").Value()
    $row3[$col] = $ws.Range($col + "3").Value()
    $row4[$col] = $ws.Range($col + "4").Value()
}

# Rotate: new row2 = old row4, new row3 = old row2, new row4 = old row3
foreach ($col in $cols) {
    $ws.Range($col + "2").Value = $row4[$col]
    $ws.Range($col + "3").Value = $row2[$col]
    $ws.Range($col + "4").Value = $row3[$col]
}
